$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 2 (shifts all rows below up by one), matching the
# "menu item deletion when branch is closed" commit.
$ws.Rows.Item(2).Delete()
